$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.5
$ws.Range("I2").Value = 1.73
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 2.38
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.88
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 23
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 41
$ws.Range("AC2").Value = 10
$ws.Range("AD2").Value = 7
$ws.Range("AG2").Value = 251
$ws.Range("AK2").Value = 13
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 26
$ws.Range("AX2").Value = 9
# Row 3
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 6.25
$ws.Range("J3").Value = 1.62
$ws.Range("K3").Value = 2.88
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5
$ws.Range("Q3").Value = 1.48
$ws.Range("R3").Value = 2.6
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("W3").Value = 9
$ws.Range("X3").Value = 7
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 7.5
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 17
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 29
$ws.Range("AL3").Value = 81
$ws.Range("AM3").Value = 67
$ws.Range("AN3").Value = 3.25
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 13
$ws.Range("AR3").Value = 34
$ws.Range("AT3").Value = 3.75
$ws.Range("AU3").Value = 10
$ws.Range("AW3").Value = 11
$ws.Range("BA3").Value = 201
$ws.Range("BB3").Value = 301
# Row 4
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.63
$ws.Range("K4").Value = 2.1
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 9.5
$ws.Range("Y4").Value = 9
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 6.5
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 11
$ws.Range("AJ4").Value = 13
$ws.Range("AO4").Value = 11
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 41
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 2.75
$ws.Range("AU4").Value = 8
$ws.Range("AY4").Value = 29
$ws.Range("BB4").Value = 201
# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 3.75
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("X5").Value = 11
$ws.Range("AC5").Value = 8.5
$ws.Range("AG5").Value = 301
$ws.Range("AH5").Value = 9
$ws.Range("AL5").Value = 26
$ws.Range("AT5").Value = 2.63
$ws.Range("AV5").Value = 51
$ws.Range("AZ5").Value = 51
$ws.Range("BB5").Value = 201
# Row 7
$ws.Range("BD7").Value = 126
# Row 8
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 3.5
$ws.Range("J8").Value = 2.75
$ws.Range("K8").Value = 2.05
$ws.Range("L8").Value = 4.33
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("X8").Value = 9
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 6.5
$ws.Range("AG8").Value = 351
$ws.Range("AH8").Value = 9.5
$ws.Range("AL8").Value = 29
$ws.Range("AR8").Value = 67
$ws.Range("AX8").Value = 21
$ws.Range("AY8").Value = 29
$ws.Range("AZ8").Value = 67
$ws.Range("BD8").Value = 151
# Row 9
$ws.Range("G9").Value = 2.2
$ws.Range("I9").Value = 3.25
$ws.Range("J9").Value = 2.88
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 4
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83
$ws.Range("X9").Value = 10
$ws.Range("AG9").Value = 301
$ws.Range("AH9").Value = 9
$ws.Range("AL9").Value = 29
$ws.Range("AZ9").Value = 67
# Row 11
$ws.Range("G11").Value = 1.62
$ws.Range("I11").Value = 5.4
$ws.Range("J11").Value = 2.15
$ws.Range("L11").Value = 5.3
$ws.Range("N11").Value = 6.9
$ws.Range("O11").Value = 1.31
$ws.Range("P11").Value = 2.87
$ws.Range("W11").Value = 5.9
$ws.Range("X11").Value = 6.9
$ws.Range("Z11").Value = 11.75
$ws.Range("AC11").Value = 9
$ws.Range("AD11").Value = 6.9
$ws.Range("AH11").Value = 13.5
$ws.Range("AI11").Value = 32
$ws.Range("AJ11").Value = 17
$ws.Range("AK11").Value = 110
$ws.Range("AL11").Value = 60
$ws.Range("AO11").Value = 7.7
$ws.Range("AP11").Value = 17
$ws.Range("AQ11").Value = 25
$ws.Range("AR11").Value = 55
$ws.Range("AT11").Value = 2.75
$ws.Range("AU11").Value = 7.4
$ws.Range("AV11").Value = 65
$ws.Range("AW11").Value = 6.9
$ws.Range("AX11").Value = 30
$ws.Range("BB11").Value = 400
# Row 12
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 2.15
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("T12").Value = 2.55
$ws.Range("V12").Value = 1.85
$ws.Range("W12").Value = 9.25
$ws.Range("AE12").Value = 14.5
$ws.Range("AH12").Value = 7.1
$ws.Range("AI12").Value = 10.25
$ws.Range("AO12").Value = 18.5
$ws.Range("AT12").Value = 2.52
$ws.Range("AU12").Value = 6.8
$ws.Range("AV12").Value = 60
$ws.Range("AW12").Value = 4
$ws.Range("AY12").Value = 18.5
$ws.Range("BB12").Value = 250
